# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals.
# The "K" column (column G) values are recalculated/rewritten with new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$updates = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 0
    6  = 0
    7  = 2
    8  = 3
    9  = 3
    10 = 0
    11 = 3
    12 = 2
    13 = 2
    14 = 1
    15 = 3
    16 = 0
    17 = 3
    18 = 2
    19 = 0
    20 = 1
    21 = 0
    22 = 2
    23 = 1
    24 = 2
    25 = 3
    26 = 0
    27 = 2
    28 = 0
    29 = 2
    30 = 1
    31 = 3
    32 = 2
    33 = 3
    34 = 2
    35 = 1
    36 = 0
    37 = 3
    38 = 1
    39 = 1
    40 = 2
    41 = 0
    42 = 1
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 1
    48 = 3
    49 = 3
    50 = 1
    51 = 2
    52 = 3
    53 = 1
    54 = 0
    55 = 0
    56 = 1
    57 = 0
    58 = 1
    59 = 1
    61 = 1
    62 = 1
    63 = 1
    64 = 3
    65 = 0
    66 = 1
    67 = 3
    68 = 1
    69 = 2
    70 = 1
    71 = 2
    72 = 1
    73 = 1
    74 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
